# Architecture diagram update for the "dev" environment slide.
#
# Summary of the change (see commit message: "changed architecture
# diagrams for development env"):
#   1. The "Group 1" AWS-Cloud-boundary group (cNvPr id=20) gets
#      shorter (its height shrinks).
#   2. The AWS logo picture nested inside that group (cNvPr id=22,
#      "Graphic 55") is resized to match the new group framing.
#   3. The lone picture "図 29" (id=30) is repositioned (moved up/left).
#   4. The lone picture "図 48" (id=49) is repositioned (moved up/left).
#   5. The elbow connector "コネクタ: カギ線 50" (id=51) that is glued to
#      that picture is resized/repositioned to match.
#   6. Five trailing shapes that duplicated the Route 53 / CloudFront /
#      S3-bucket icon groups and a stray picture+connector pair are
#      removed entirely (cNvPr ids 64, 67, 70, 74, 75).
#
# Implementation notes on this COM host's quirks (reverse-engineered by
# probing against the loaded deck):
#
#  * Shape.Left/Top/Width/Height setters convert the incoming point
#    value to EMU via a straight value/72*914400 and then truncate
#    toward zero, rather than rounding to the nearest EMU, before they
#    land in the OOXML <a:off>/<a:ext>. Feeding the EMU-exact point
#    value can therefore land one EMU short of the intended target.
#    Adding a nudge of +0.5 EMU worth of points before conversion fixes
#    this without perturbing already-exact values.
#  * For a shape that lives inside a group (GroupItems), the
#    Height/Width/Left/Top setters write straight into that shape's own
#    (group-local / child) coordinate space -- they do NOT re-apply the
#    parent group's off/ext -> chOff/chExt scale factor on write (even
#    though *reading* those same properties does apply the scale and
#    returns slide-space values). So, to hit a specific child-space EMU
#    value in the XML, the nested shape must be fed the point
#    equivalent of that child-space EMU value directly.
#
# Shapes are looked up by their (stable) cNvPr id rather than by a
# fixed Shapes.Item(N) index, so this script does not depend on
# collection ordering.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU-per-point, plus half an EMU's worth of points so the host's
# truncating pt->EMU conversion rounds to the exact target EMU value.
$EmuPerPt = 914400.0 / 72.0
$HalfEmuPt = 0.5 / $EmuPerPt

function EmuPt([double]$emu) {
    return ($emu / $EmuPerPt) + $HalfEmuPt
}

function Get-ShapeById($shapes, [int]$id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $candidate = $shapes.Item($i)
        if ($candidate.Id -eq $id) {
            return $candidate
        }
    }
    return $null
}

# --- 1: "Group 1" (id=20) shrinks vertically (5057192 -> 4423825 EMU). ---
$group1 = Get-ShapeById $s.Shapes 20
$group1.Height = EmuPt 4423825

# --- 2: nested AWS-logo picture "Graphic 55" (id=22) is resized (in
#     its own group-local coordinates) to 93099 EMU tall. ---
$awsLogoPic = Get-ShapeById $group1.GroupItems 22
$awsLogoPic.Height = EmuPt 93099

# --- 3: picture "図 29" (id=30) moves up/left. ---
$pic29 = Get-ShapeById $s.Shapes 30
$pic29.Left = EmuPt 911185
$pic29.Top = EmuPt 2405023

# --- 4: picture "図 48" (id=49) moves up/left. ---
$pic48 = Get-ShapeById $s.Shapes 49
$pic48.Left = EmuPt 6747224
$pic48.Top = EmuPt 5099883

# --- 5: elbow connector "コネクタ: カギ線 50" (id=51) shrinks to match
#     the picture it is anchored to. ---
$connector50 = Get-ShapeById $s.Shapes 51
$connector50.Left = EmuPt 5956094
$connector50.Top = EmuPt 3809937
$connector50.Width = EmuPt 791130
$connector50.Height = EmuPt 1565366

# --- 6: remove the trailing Route 53 / CloudFront / S3-bucket groups
#     plus the stray picture+connector pair that duplicated the
#     frontend-deployment icon. Delete from the highest index down so
#     earlier indices stay valid while iterating. ---
$idsToDelete = @(64, 67, 70, 74, 75)
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shape = $s.Shapes.Item($i)
    if ($idsToDelete -contains $shape.Id) {
        $shape.Delete()
    }
}
